$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2, D2, E2 deleted (cleared); C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -3.7869223958651617
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: B3:E3 updated
$ws.Range("B3").Value = -3.7209366071500796
$ws.Range("C3").Value = 1.8519109904298858
$ws.Range("D3").Value = -3.2448504095349477
$ws.Range("E3").Value = 11.02810119935404

# Update selection to B1:E3
$ws.Range("B1:E3").Select() | Out-Null
